$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.690.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.20%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.264.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.89%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'304.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.09%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'91.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.81%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.529"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.60%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +1.15%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'32.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.62%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +0.08%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.47%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.75%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.18%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.615.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.18%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'14.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.00%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.243.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.04%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.764"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.14%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'41.605.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.31%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'12.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +9.15%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.0₃0903"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.30%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +1.82%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'66.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.26%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'240.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.58%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +3.64%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.14%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +4.80%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'23.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.27%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +1.30%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -4.45%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'160.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.65%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'34.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.78%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.02%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.18%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0743"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.39%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.70%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.84%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +2.18%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +1.84%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'16.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.42%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.81%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'3.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.96%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.047.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.10%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'19.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.13%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'10.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.18%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.0278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.26%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.24%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +7.48%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.53%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'72.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.13%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.37%  "
$ws.Range("E51").Style = "Normal"
